$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 20:55"

# Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 4942293
$ws.Range("C4").Value = 23873
$ws.Range("D4").Value = 2503460
$ws.Range("E4").Value = 2277835
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 708
$ws.Range("H4").Value = 160998

# India (row 6) - updated case counts
$ws.Range("B6").Value = 1963239
$ws.Range("C6").Value = 56626
$ws.Range("D6").Value = 1327200
$ws.Range("E6").Value = 595300
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 919
$ws.Range("H6").Value = 40739

# Francia (row 23) - updated case counts
$ws.Range("B23").Value = 194029
$ws.Range("C23").Value = 1695
$ws.Range("D23").Value = 82166
$ws.Range("E23").Value = 81558
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 30305

# Catar (row 28) - updated case counts
$ws.Range("B28").Value = 111805
$ws.Range("C28").Value = 267
$ws.Range("D28").Value = 108539
$ws.Range("E28").Value = 3088
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 178

# Libano overtakes Croacia and Republica de Yibuti in the ranking (rows 100-102)
$ws.Range("A100").Value = "Libano"
$ws.Range("B100").Value = 5417
$ws.Range("C100").Value = 146
$ws.Range("D100").Value = 1880
$ws.Range("E100").Value = 3469
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 68

$ws.Range("A101").Value = "Croacia"
$ws.Range("B101").Value = 5376
$ws.Range("C101").Value = 58
$ws.Range("D101").Value = 4589
$ws.Range("E101").Value = 633
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 154

$ws.Range("A102").Value = "Republica de Yibuti"
$ws.Range("B102").Value = 5330
$ws.Range("C102").Value = 82
$ws.Range("D102").Value = 5057
$ws.Range("E102").Value = 214
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 59

# Malaui (row 108) - updated case counts
$ws.Range("B108").Value = 4426
$ws.Range("C108").Value = 65
$ws.Range("D108").Value = 2078
$ws.Range("E108").Value = 2212
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 8
$ws.Range("H108").Value = 136

# Mayotte (row 117) - updated case counts
$ws.Range("B117").Value = 3031
$ws.Range("C117").Value = 8
$ws.Range("D117").Value = 2738
$ws.Range("E117").Value = 254
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 39

# Suazilandia (row 118) - updated case counts
$ws.Range("B118").Value = 2909
$ws.Range("C118").Value = 53
$ws.Range("D118").Value = 1385
$ws.Range("E118").Value = 1471
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 4
$ws.Range("H118").Value = 53

# Sri Lanka (row 119) - updated case counts
$ws.Range("B119").Value = 2839
$ws.Range("C119").Value = 5
$ws.Range("D119").Value = 2537
$ws.Range("E119").Value = 291
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 11

# Cabo Verde (row 121) - updated case counts
$ws.Range("B121").Value = 2689
$ws.Range("C121").Value = 58
$ws.Range("D121").Value = 1955
$ws.Range("E121").Value = 707
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 27

# Mali (row 122) - updated case counts
$ws.Range("B122").Value = 2546
$ws.Range("C122").Value = 3
$ws.Range("D122").Value = 1950
$ws.Range("E122").Value = 472
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 124

# Timor Oriental overtakes Santa Lucia in the ranking (rows 202-203)
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("B202").Value = 25
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 24
$ws.Range("E202").Value = 1
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 25
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 24
$ws.Range("E203").Value = 1
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0
